# Applies the workbook edit described by the diff:
#  - Updates the numeric value + formatting of the "4" cell on each of the
#    three sheets to 4.269 with a "0.00" number format.
#  - Updates each sheet's saved cursor/selection to the new active cell.
#  - Leaves "Tableless" as the active (frontmost) sheet/tab.

$wb = $excel.ActiveWorkbook

# --- WithTable sheet: A5 4 -> 4.269, selection E23 -> A5 ---
$wsWithTable = $wb.Worksheets.Item("WithTable")
$wsWithTable.Range("A5").Value = 4.269
$wsWithTable.Range("A5").NumberFormat = "0.00"
$wsWithTable.Activate()
$wsWithTable.Range("A5").Select() | Out-Null

# --- WithTable_Duplicate sheet: B8 4 -> 4.269, selection M11 -> B8 ---
$wsDuplicate = $wb.Worksheets.Item("WithTable_Duplicate")
$wsDuplicate.Range("B8").Value = 4.269
$wsDuplicate.Range("B8").NumberFormat = "0.00"
$wsDuplicate.Activate()
$wsDuplicate.Range("B8").Select() | Out-Null

# --- Tableless sheet: A5 4 -> 4.269, selection I19 -> A5 ---
$wsTableless = $wb.Worksheets.Item("Tableless")
$wsTableless.Range("A5").Value = 4.269
$wsTableless.Range("A5").NumberFormat = "0.00"
$wsTableless.Activate()
$wsTableless.Range("A5").Select() | Out-Null

# Tableless ends up as the active tab, matching the target workbook view.
$wsTableless.Activate()
